$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert 4 new "Snacks" rows above the old "Kinder Creamy" row
#    first (so every subsequent write lands on its final row number),
#    then fill in the brand-new text in the particular sequence that
#    reproduces the shared-string build order Excel itself produced
#    while authoring the change.
# ------------------------------------------------------------------
$ws.Rows.Item(49).Insert() | Out-Null
$ws.Rows.Item(49).Insert() | Out-Null
$ws.Rows.Item(49).Insert() | Out-Null
$ws.Rows.Item(49).Insert() | Out-Null

$ws.Cells.Item(69, 2).Value = "Hell"
$ws.Cells.Item(69, 5).Value = "Hell.jpg"
$ws.Cells.Item(70, 2).Value = "Big Bottel Cold Drink"
$ws.Cells.Item(70, 5).Value = "Big Bottle 85.jpg"

$ws.Cells.Item(50, 2).Value = "Popz"
$ws.Cells.Item(52, 2).Value = "Duo"
$ws.Cells.Item(49, 2).Value = "Act II Sweet And Salty"
$ws.Cells.Item(49, 5).Value = "Act II Sweet And Salty.jpg"
$ws.Cells.Item(50, 5).Value = "Popz.jpg"
$ws.Cells.Item(51, 2).Value = "Act II Pop corn Butter"
$ws.Cells.Item(52, 5).Value = "Duo 5.jpg"
$ws.Cells.Item(51, 5).Value = "Act II Pop corn Butter.jpg"

# "Thumsup25" (Rs. 25) becomes "Thumsup30" (Rs. 30); after the insert
# above this item now lives on row 61.
$ws.Cells.Item(61, 2).Value = "Thumsup30"

# ------------------------------------------------------------------
# 2) Fill in the remaining (already-existing-shared-string / numeric)
#    columns for every new/changed row.
# ------------------------------------------------------------------
$ws.Cells.Item(49, 1).Value = "Snacks"
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 25
$ws.Cells.Item(49, 6).Value = "Fast Food"

$ws.Cells.Item(50, 1).Value = "Snacks"
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 10
$ws.Cells.Item(50, 6).Value = "Fast Food"

$ws.Cells.Item(51, 1).Value = "Snacks"
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 10
$ws.Cells.Item(51, 6).Value = "Fast Food"

$ws.Cells.Item(52, 1).Value = "Snacks"
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 5
$ws.Cells.Item(52, 6).Value = "Fast Food"

$ws.Cells.Item(61, 4).Value = 30

$ws.Cells.Item(69, 1).Value = "Drinks"
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 60
$ws.Cells.Item(69, 6).Value = "Fast Food"

$ws.Cells.Item(70, 1).Value = "Drinks"
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(70, 4).Value = 85
$ws.Cells.Item(70, 6).Value = "Fast Food"

# ------------------------------------------------------------------
# 3) The hidden _FilterDatabase name tracked the old table extent
#    (…$F$61); widen it to the new last row (…$F$65).
# ------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "_xlnm._FilterDatabase" -or $n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$65"
    }
}

# ------------------------------------------------------------------
# 4) Refresh the view: scroll/selection moved to B62 with top row 48.
# ------------------------------------------------------------------
$ws.Range("A48").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 48
$ws.Range("B62").Select() | Out-Null
